$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.192
$ws.Range("B3").Value = 5.933
$ws.Range("E3").Value = 16.263
$ws.Range("E12").Value = 17.924
$ws.Range("B14").Value = 5.756
$ws.Range("B21").Value = 9.238000000000001
$ws.Range("B23").Value = 7.179
$ws.Range("E24").Value = 16.842
$ws.Range("B25").Value = 5.029
$ws.Range("D25").Value = -7.648000000000001
$ws.Range("E25").Value = 17.157
$ws.Range("B26").Value = 6.182
$ws.Range("D27").Value = -8.5
$ws.Range("B29").Value = 6.006
$ws.Range("D31").Value = -8.222
$ws.Range("D39").Value = -7.782999999999999
$ws.Range("D48").Value = -7.414999999999999
$ws.Range("E50").Value = 16.45
$ws.Range("D51").Value = -8.294
$ws.Range("D52").Value = -8.054
$ws.Range("B53").Value = 6.455
$ws.Range("E53").Value = 16.82
$ws.Range("D55").Value = -8.035
$ws.Range("D56").Value = -8.135999999999999
$ws.Range("B57").Value = 4.925
$ws.Range("D57").Value = -8.059999999999999
$ws.Range("E57").Value = 16.623
$ws.Range("B59").Value = 5.145
$ws.Range("E61").Value = 16.628
$ws.Range("E63").Value = 17.656
$ws.Range("B69").Value = 5.08
$ws.Range("E70").Value = 17.594
$ws.Range("D73").Value = -8.278000000000002
$ws.Range("B79").Value = 6.266
$ws.Range("B83").Value = 5.915
$ws.Range("E86").Value = 16.394
$ws.Range("D89").Value = -6.356
$ws.Range("D90").Value = -7.507
$ws.Range("B91").Value = 5.600000000000001
$ws.Range("D92").Value = -6.62
$ws.Range("B93").Value = 5.898999999999999
$ws.Range("E98").Value = 16.487
$ws.Range("E100").Value = 16.617
$ws.Range("E102").Value = 16.205
